$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'277.02"
$ws.Range("E2").Value = "'0.98%"
$ws.Range("D3").Value = "'27.13"
$ws.Range("E3").Value = "'1.29%"
$ws.Range("D4").Value = "'4.848"
$ws.Range("E4").Value = "'-0.07%"
$ws.Range("D5").Value = "'0.06405"
$ws.Range("E5").Value = "'1.28%"
$ws.Range("D6").Value = "'6.935"
$ws.Range("E6").Value = "'0.46%"
$ws.Range("D7").Value = "'1.200"
$ws.Range("E7").Value = "'-5.75%"
$ws.Range("D8").Value = "'0.8760"
$ws.Range("E8").Value = "'0.47%"
$ws.Range("D9").Value = "'0.1523"
$ws.Range("E9").Value = "'4.27%"
$ws.Range("D10").Value = "'0.05123"
$ws.Range("E10").Value = "'2.53%"
$ws.Range("D11").Value = "'0.07495"
$ws.Range("E11").Value = "'2.34%"
$ws.Range("D12").Value = "'0.02966"
$ws.Range("E12").Value = "'-0.18%"
$ws.Range("D13").Value = "'0.08979"
$ws.Range("E13").Value = "'-0.58%"
$ws.Range("D14").Value = "'0.001560"
$ws.Range("E14").Value = "'-0.63%"
$ws.Range("D15").Value = "'0.0006343"
$ws.Range("E15").Value = "'0.45%"
$ws.Range("D16").Value = "'0.006073"
$ws.Range("E16").Value = "'0.67%"
$ws.Range("D17").Value = "'3.477"
$ws.Range("E17").Value = "'0.82%"
$ws.Range("D18").Value = "'3.307"
$ws.Range("E18").Value = "'-0.57%"
$ws.Range("E19").Value = "'-0.40%"
$ws.Range("E21").Value = "'1.04%"
$ws.Range("D22").Value = "'3.914"
$ws.Range("E22").Value = "'0.21%"
$ws.Range("D23").Value = "'0.04425"
$ws.Range("E23").Value = "'1.77%"
$ws.Range("D25").Value = "'0.001175"
$ws.Range("E25").Value = "'-0.25%"
$ws.Range("D26").Value = "'0.003865"
$ws.Range("E26").Value = "'-9.09%"
$ws.Range("E27").Value = "'8.20%"
$ws.Range("E28").Value = "'15.00%"
$ws.Range("D40").Value = "'0.04159"
$ws.Range("E40").Value = "'2.66%"
$ws.Range("D41").Value = "'0.006806"
$ws.Range("E41").Value = "'1.39%"
$ws.Range("E42").Value = "'0.51%"
$ws.Range("D43").Value = "'0.001949"
$ws.Range("E43").Value = "'-7.28%"
$ws.Range("D44").Value = "'0.01193"
$ws.Range("E44").Value = "'11.46%"
$ws.Range("D45").Value = "'0.00005303"
$ws.Range("E45").Value = "'-0.15%"
$ws.Range("E46").Value = "'13.49%"
$ws.Range("D47").Value = "'0.01851"
$ws.Range("E47").Value = "'-7.52%"
